# Agregar nueva metrica "Promedio arriendo mes, semana, año." al excel de Progress.
# Se inserta como nueva fila 27 (categoria "Por País"), desplazando hacia abajo las
# filas existentes "Usuario con reviews en más países" y
# "Comparativa de precio promedio entre países".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Desplazar el contenido de las filas 27-28 una posicion hacia abajo (28->29, 27->28)
# copiando solo el rango de datos (A:D) para no arrastrar formato de columnas vacias.
$ws.Range("A28:D28").Copy($ws.Range("A29:D29"))
$ws.Range("A27:D27").Copy($ws.Range("A28:D28"))

# Completar la nueva fila 27 con la metrica agregada.
$ws.Range("A27").Value = "Promedio arriendo mes, semana, año."
$ws.Range("B27").Value = 0
$ws.Range("D27").Value = "Por País"

# Actualizar el rango de formato condicional (columna B de estado) para cubrir la fila nueva.
$fcs = $ws.Range("B7:B28").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("B7:B29"))
}
# Asegurar que el formato diferencial (dxf) de las reglas quede definido en el libro.
$fcs.Item(1).Font.Size = 10

# Ajustar levemente el color de fondo del indicador "Finished" de la leyenda.
$ws.Range("G6").Interior.PatternColor = 32768

# Mover la celda activa/seleccion a D28, como quedo tras editar la planilla.
$ws.Range("D28").Select()
